# Generate Report for Archive
# The handoff status text moves from "Ready for handoff" to "In Translation"
# everywhere it is reported: the per-locale status column on the Overview
# sheet (zh-cn / de-de columns) and the Status column on each locale sheet.
# Excel then re-auto-sizes those status columns to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns("E:F").AutoFit()

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns("C:C").AutoFit()

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns("C:C").AutoFit()
